# Automatische test-sync: 2025-06-22 18:58:50
# Appends a new "Onjuiste bestelling" log entry to the Logs sheet and
# refreshes the Dashboard's category-count table to reflect it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: add the new row 28
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A28").Value = "Onjuiste bestelling"
$logs.Range("B28").Value = "mailmind.test@zohomail.eu"
$logs.Range("C28").Value = "Ik heb een ander product ontvangen dan ik had besteld."
$logs.Range("D28").Value = "Retour / Terugbetaling"
$logs.Range("E28").Value = "Beste klant,
Dank u voor uw bericht. We vinden het vervelend om te horen dat u een ander product heeft ontvangen dan u had besteld. Om dit probleem op te lossen, willen we u vragen om het volgende te doen:
1. Stuur ons alstublieft een foto van het ontvangen product, samen met uw bestelnummer, naar ons e-mailadres, zodat we dit kunnen verifiëren.
2. Geef ons ook de omschrijving van het product dat u had besteld, zodat we het verschil kunnen vaststellen.
Zodra we deze informatie hebben ontvangen, zullen we ons best doen om dit zo snel mogelijk voor u op te lossen.
Met vriendelijke groet,
[Bedrijfsnaam] E-mailassistent"
$logs.Range("F28").Value = "2025-06-22 18:58:15"
$logs.Range("G28").Value = "Ja"

# Extend the conditional-formatting ranges (D2:D27 -> D2:D28, G2:G27 -> G2:G28).
# Modifying one rule's AppliesTo range updates the shared sqref for every
# cfRule that lives inside the same <conditionalFormatting> block.
$dRules = $logs.Range("D2:D27").FormatConditions
$dRules.Item(1).ModifyAppliesToRange($logs.Range("D2:D28"))

$gRules = $logs.Range("G2:G27").FormatConditions
$gRules.Item(1).ModifyAppliesToRange($logs.Range("G2:G28"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: re-rank the category counts
#    ("Retour / Terugbetaling" now has 3 hits, "Productinformatie" 3,
#    "Samenwerking / Partnerverzoek" stays at 2 but drops further down).
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Retour / Terugbetaling"
$dash.Range("B4").Value = 3

$dash.Range("A5").Value = "Productinformatie"
$dash.Range("B5").Value = 3

$dash.Range("A9").Value = "Samenwerking / Partnerverzoek"
$dash.Range("B9").Value = 2
